$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: the "black gate" quest entry loses its Quest1/2/3 sub-columns and
# gets the new closing-story text in the Descript column.
$ws.Range("E9").Value = "|通过三把钥匙，我打开了高塔的大门，一切的谜题感觉就要在今天解开了。那我传说中的的勇士在塔内究竟找到了什么？高塔第一层都是一些关于塔历史的壁画，似乎没有什么奇怪的东西。找到阶梯后，我就进入到了第二层，这里充满了各种机关，并且从机关的类型看来，年代并不是非常久远。我开始怀疑这个塔从何而来，这里并不像一个古代文明的古迹。奇怪的是，越是接近塔顶，这里黑暗的力量就越强大。终于，我来到了塔顶，在这里发现了一扇异世界的传送门，而所有黑暗力量的怪物都是由此而来。再消灭了从门内出现的一批怪物后，我赶紧破坏掉了传送门。一切都结束了吗？"
$ws.Range("F9:K9").ClearContents()

# Row 10 becomes the new "newbie quest" (12001001 / beach), replacing the
# old 12000007 row.
$ws.Range("A10").Value = 12001001
$ws.Range("B10").Value = "海滩"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = "铁匠和我说了很多关于我叔叔以前的故事，原来他曾经是一个十分厉害的卡片召唤师。"

# Drop the now-finished quest rows 11-18.
$ws.Range("A11:K18").EntireRow.Delete()

$ws.Range("D10").Select()
